$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reassign reviewers on row 2 (story 1.1) ---
$ws.Range("D2").Value = "Ahmed Ashraf"
$ws.Range("E2").Value = "Abed Hossam"
$ws.Range("F2").Value = "Omar Bakr"

# --- Reassign reviewer on row 5 (story 2.1) ---
$ws.Range("D5").Value = "Mostafa Waleed"

# --- Insert a new backlog item (story 3.4) after row 9, before the blank separator row ---
$ws.Rows("10").Insert()

$ws.Range("A10").Value = 3.4
$ws.Range("B10").Value = "As a user, I want to view my notifications"
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = "Ahmed Ashraf"
$ws.Range("E10").Value = "Omar Khaled"
$ws.Range("F10").Value = "Yosri Khaled"

$ws.Range("B9").Select()
